$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells(54, 4).Value = 44434
$ws.Cells(55, 4).Value = 44340
$ws.Cells(55, 10).Value = 700
$ws.Cells(55, 11).Value = 5000
$ws.Cells(55, 12).Value = 5500
$ws.Cells(55, 13).Value = 5250
$ws.Cells(55, 16).Value = 262
$ws.Cells(56, 4).Value = 44174
$ws.Cells(56, 10).Value = 720
$ws.Cells(56, 11).Value = 4500
$ws.Cells(56, 12).Value = 5000
$ws.Cells(56, 13).Value = 4750
$ws.Cells(56, 16).Value = 238
$ws.Cells(57, 4).Value = 44221
$ws.Cells(57, 10).Value = 520
$ws.Cells(57, 11).Value = 5000
$ws.Cells(57, 12).Value = 5500
$ws.Cells(57, 13).Value = 5250
$ws.Cells(57, 16).Value = 262
$ws.Cells(58, 4).Value = 44271
$ws.Cells(58, 10).Value = 600
$ws.Cells(58, 11).Value = 5500
$ws.Cells(58, 12).Value = 6000
$ws.Cells(58, 13).Value = 5750
$ws.Cells(58, 16).Value = 288
$ws.Cells(59, 4).Value = 44280
$ws.Cells(59, 10).Value = 680
$ws.Cells(60, 4).Value = 44299
$ws.Cells(60, 10).Value = 660
$ws.Cells(60, 11).Value = 5000
$ws.Cells(60, 12).Value = 5500
$ws.Cells(60, 13).Value = 5250
$ws.Cells(60, 16).Value = 262
$ws.Cells(61, 4).Value = 44428
$ws.Cells(61, 10).Value = 720
$ws.Cells(61, 11).Value = 5000
$ws.Cells(61, 12).Value = 5500
$ws.Cells(61, 13).Value = 5250
$ws.Cells(61, 16).Value = 262
$ws.Cells(62, 4).Value = 44162
$ws.Cells(62, 10).Value = 540
$ws.Cells(62, 11).Value = 4800
$ws.Cells(62, 12).Value = 5000
$ws.Cells(62, 13).Value = 4900
$ws.Cells(62, 16).Value = 245
$ws.Cells(63, 4).Value = 44342
$ws.Cells(63, 10).Value = 800
$ws.Cells(63, 11).Value = 4800
$ws.Cells(63, 12).Value = 5000
$ws.Cells(63, 13).Value = 4900
$ws.Cells(63, 16).Value = 245
$ws.Cells(64, 4).Value = 44362
$ws.Cells(64, 10).Value = 700
$ws.Cells(65, 4).Value = 44175
$ws.Cells(65, 10).Value = 600
$ws.Cells(65, 11).Value = 5000
$ws.Cells(65, 12).Value = 5500
$ws.Cells(65, 13).Value = 5250
$ws.Cells(65, 16).Value = 262
$ws.Cells(66, 4).Value = 44384
$ws.Cells(66, 10).Value = 800
$ws.Cells(67, 4).Value = 44242
$ws.Cells(67, 10).Value = 680
$ws.Cells(67, 11).Value = 5500
$ws.Cells(67, 12).Value = 6000
$ws.Cells(67, 13).Value = 5750
$ws.Cells(67, 16).Value = 288
$ws.Cells(68, 4).Value = 44204
$ws.Cells(68, 10).Value = 740
$ws.Cells(69, 4).Value = 44323
$ws.Cells(69, 11).Value = 4800
$ws.Cells(69, 12).Value = 5000
$ws.Cells(69, 13).Value = 4900
$ws.Cells(69, 16).Value = 245
$ws.Cells(70, 4).Value = 44200
$ws.Cells(70, 10).Value = 520
$ws.Cells(71, 4).Value = 44363
$ws.Cells(71, 11).Value = 5000
$ws.Cells(71, 12).Value = 5500
$ws.Cells(71, 13).Value = 5250
$ws.Cells(71, 16).Value = 262
$ws.Cells(72, 4).Value = 44216
$ws.Cells(73, 4).Value = 44349
$ws.Cells(73, 11).Value = 4800
$ws.Cells(73, 12).Value = 5000
$ws.Cells(73, 13).Value = 4900
$ws.Cells(73, 16).Value = 245
$ws.Cells(74, 4).Value = 44385
$ws.Cells(74, 10).Value = 720
$ws.Cells(74, 11).Value = 5000
$ws.Cells(74, 12).Value = 5500
$ws.Cells(74, 13).Value = 5250
$ws.Cells(74, 16).Value = 262
$ws.Cells(75, 4).Value = 44258
$ws.Cells(75, 10).Value = 800
$ws.Cells(75, 11).Value = 5500
$ws.Cells(75, 12).Value = 6000
$ws.Cells(75, 13).Value = 5750
$ws.Cells(75, 16).Value = 288
$ws.Cells(76, 4).Value = 44243
$ws.Cells(76, 10).Value = 600
$ws.Cells(76, 11).Value = 5500
$ws.Cells(76, 12).Value = 6000
$ws.Cells(76, 13).Value = 5750
$ws.Cells(76, 16).Value = 288
$ws.Cells(77, 4).Value = 44427
$ws.Cells(77, 10).Value = 700
$ws.Cells(77, 11).Value = 5000
$ws.Cells(77, 12).Value = 5500
$ws.Cells(77, 13).Value = 5250
$ws.Cells(77, 16).Value = 262
$ws.Cells(78, 4).Value = 44413
$ws.Cells(78, 10).Value = 720
$ws.Cells(79, 4).Value = 44176
$ws.Cells(79, 10).Value = 560
$ws.Cells(79, 11).Value = 4800
$ws.Cells(79, 12).Value = 5000
$ws.Cells(79, 13).Value = 4900
$ws.Cells(79, 16).Value = 245
$ws.Cells(80, 4).Value = 44421
$ws.Cells(80, 10).Value = 700
$ws.Cells(80, 11).Value = 5000
$ws.Cells(80, 12).Value = 5500
$ws.Cells(80, 13).Value = 5250
$ws.Cells(80, 16).Value = 262
$ws.Cells(81, 4).Value = 44222
$ws.Cells(81, 10).Value = 600
$ws.Cells(81, 11).Value = 5000
$ws.Cells(81, 12).Value = 5500
$ws.Cells(81, 13).Value = 5250
$ws.Cells(81, 16).Value = 262
$ws.Cells(82, 4).Value = 44237
$ws.Cells(82, 11).Value = 5500
$ws.Cells(82, 12).Value = 6000
$ws.Cells(82, 13).Value = 5750
$ws.Cells(82, 16).Value = 288
$ws.Cells(83, 4).Value = 44273
$ws.Cells(84, 4).Value = 44377
$ws.Cells(84, 10).Value = 800
$ws.Cells(85, 4).Value = 44257
$ws.Cells(85, 10).Value = 700
$ws.Cells(85, 11).Value = 6000
$ws.Cells(85, 12).Value = 6500
$ws.Cells(85, 13).Value = 6250
$ws.Cells(85, 16).Value = 312
$ws.Cells(86, 4).Value = 44400
$ws.Cells(86, 10).Value = 720
$ws.Cells(87, 4).Value = 44426
$ws.Cells(88, 4).Value = 44225
$ws.Cells(88, 10).Value = 760
$ws.Cells(89, 4).Value = 44295
$ws.Cells(89, 10).Value = 800
$ws.Cells(90, 4).Value = 44194
$ws.Cells(90, 10).Value = 560
$ws.Cells(91, 4).Value = 44390
$ws.Cells(91, 10).Value = 600
$ws.Cells(92, 4).Value = 44321
$ws.Cells(92, 10).Value = 800
$ws.Cells(93, 4).Value = 44298
$ws.Cells(93, 10).Value = 680
$ws.Cells(94, 4).Value = 44383
$ws.Cells(94, 10).Value = 600
$ws.Cells(95, 4).Value = 44354
$ws.Cells(95, 10).Value = 700
$ws.Cells(96, 4).Value = 44169
$ws.Cells(96, 10).Value = 540
$ws.Cells(97, 4).Value = 44410
$ws.Cells(97, 10).Value = 760
$ws.Cells(97, 11).Value = 5000
$ws.Cells(97, 12).Value = 5500
$ws.Cells(97, 13).Value = 5250
$ws.Cells(97, 16).Value = 262
$ws.Cells(98, 4).Value = 44412
$ws.Cells(99, 4).Value = 44336
$ws.Cells(99, 10).Value = 600
$ws.Cells(99, 11).Value = 4800
$ws.Cells(99, 12).Value = 5000
$ws.Cells(99, 13).Value = 4900
$ws.Cells(99, 16).Value = 245
$ws.Cells(100, 4).Value = 44300
$ws.Cells(100, 10).Value = 800
$ws.Cells(100, 11).Value = 5000
$ws.Cells(100, 12).Value = 5500
$ws.Cells(100, 13).Value = 5250
$ws.Cells(100, 16).Value = 262
$ws.Cells(101, 4).Value = 44172
$ws.Cells(101, 10).Value = 760
$ws.Cells(102, 4).Value = 44214
$ws.Cells(102, 10).Value = 540
$ws.Cells(102, 11).Value = 4500
$ws.Cells(102, 12).Value = 5000
$ws.Cells(102, 13).Value = 4750
$ws.Cells(102, 16).Value = 238
$ws.Cells(103, 4).Value = 44371
$ws.Cells(103, 10).Value = 700
$ws.Cells(103, 11).Value = 5000
$ws.Cells(103, 12).Value = 5500
$ws.Cells(103, 13).Value = 5250
$ws.Cells(103, 16).Value = 262
$ws.Cells(104, 4).Value = 44238
$ws.Cells(104, 10).Value = 680
$ws.Cells(104, 11).Value = 5500
$ws.Cells(104, 12).Value = 6000
$ws.Cells(104, 13).Value = 5750
$ws.Cells(104, 16).Value = 288
$ws.Cells(105, 4).Value = 44274
$ws.Cells(105, 11).Value = 6000
$ws.Cells(105, 12).Value = 6500
$ws.Cells(105, 13).Value = 6250
$ws.Cells(105, 16).Value = 312
$ws.Cells(106, 4).Value = 44320
$ws.Cells(106, 10).Value = 600
$ws.Cells(107, 4).Value = 44314
$ws.Cells(108, 4).Value = 44364
$ws.Cells(108, 10).Value = 680
$ws.Cells(109, 4).Value = 44435
$ws.Cells(109, 10).Value = 3276
$ws.Cells(109, 11).Value = 5000
$ws.Cells(109, 12).Value = 5500
$ws.Cells(109, 13).Value = 5275
$ws.Cells(109, 16).Value = 264
$ws.Cells(110, 4).Value = 44431
$ws.Cells(110, 10).Value = 720
$ws.Cells(110, 11).Value = 5000
$ws.Cells(110, 12).Value = 5500
$ws.Cells(110, 13).Value = 5250
$ws.Cells(110, 16).Value = 262
$ws.Cells(111, 4).Value = 44405
$ws.Cells(111, 10).Value = 800
$ws.Cells(112, 4).Value = 44224
$ws.Cells(113, 4).Value = 44260
$ws.Cells(113, 10).Value = 800
$ws.Cells(113, 11).Value = 5500
$ws.Cells(113, 12).Value = 6000
$ws.Cells(113, 13).Value = 5750
$ws.Cells(113, 16).Value = 288
$ws.Cells(114, 4).Value = 44327
$ws.Cells(114, 10).Value = 600
$ws.Cells(115, 4).Value = 44209
$ws.Cells(115, 10).Value = 720
$ws.Cells(116, 4).Value = 44231
$ws.Cells(116, 10).Value = 600
$ws.Cells(117, 4).Value = 44313
$ws.Cells(117, 10).Value = 600
$ws.Cells(117, 11).Value = 5000
$ws.Cells(117, 12).Value = 5500
$ws.Cells(117, 13).Value = 5250
$ws.Cells(117, 16).Value = 262
$ws.Cells(118, 4).Value = 44330
$ws.Cells(118, 10).Value = 800
$ws.Cells(119, 4).Value = 44391
$ws.Cells(119, 10).Value = 800
$ws.Cells(120, 4).Value = 44193
$ws.Cells(120, 10).Value = 540
$ws.Cells(121, 4).Value = 44351
$ws.Cells(121, 10).Value = 800
$ws.Cells(121, 11).Value = 4800
$ws.Cells(121, 12).Value = 5000
$ws.Cells(121, 13).Value = 4900
$ws.Cells(121, 16).Value = 245
$ws.Cells(122, 4).Value = 44350
$ws.Cells(122, 10).Value = 700
$ws.Cells(122, 11).Value = 4800
$ws.Cells(122, 12).Value = 5000
$ws.Cells(122, 13).Value = 4900
$ws.Cells(122, 16).Value = 245
$ws.Cells(123, 4).Value = 44196
$ws.Cells(123, 10).Value = 600
$ws.Cells(124, 4).Value = 44315
$ws.Cells(124, 10).Value = 680
$ws.Cells(124, 11).Value = 5000
$ws.Cells(124, 12).Value = 5500
$ws.Cells(124, 13).Value = 5250
$ws.Cells(124, 16).Value = 262
$ws.Cells(125, 4).Value = 44358
$ws.Cells(125, 10).Value = 700
$ws.Cells(125, 11).Value = 5000
$ws.Cells(125, 12).Value = 5500
$ws.Cells(125, 13).Value = 5250
$ws.Cells(125, 16).Value = 262
$ws.Cells(126, 4).Value = 44389
$ws.Cells(126, 10).Value = 760
$ws.Cells(127, 4).Value = 44399
$ws.Cells(127, 10).Value = 720
$ws.Cells(128, 4).Value = 44251
$ws.Cells(128, 10).Value = 800
$ws.Cells(128, 11).Value = 5500
$ws.Cells(128, 12).Value = 6000
$ws.Cells(128, 13).Value = 5750
$ws.Cells(128, 16).Value = 288
$ws.Cells(129, 4).Value = 44250
$ws.Cells(129, 10).Value = 600
$ws.Cells(129, 11).Value = 5500
$ws.Cells(129, 12).Value = 6000
$ws.Cells(129, 13).Value = 5750
$ws.Cells(129, 16).Value = 288
$ws.Cells(130, 4).Value = 44292
$ws.Cells(130, 10).Value = 600
$ws.Cells(131, 4).Value = 44305
$ws.Cells(131, 10).Value = 680
$ws.Cells(132, 4).Value = 44294
$ws.Cells(132, 10).Value = 600
$ws.Cells(132, 11).Value = 5000
$ws.Cells(132, 12).Value = 5500
$ws.Cells(132, 13).Value = 5250
$ws.Cells(132, 16).Value = 262
$ws.Cells(133, 4).Value = 44417
$ws.Cells(134, 4).Value = 44419
$ws.Cells(134, 10).Value = 800
$ws.Cells(135, 4).Value = 44420
$ws.Cells(135, 10).Value = 700
$ws.Cells(136, 4).Value = 44245
$ws.Cells(136, 10).Value = 700
$ws.Cells(136, 11).Value = 5500
$ws.Cells(136, 12).Value = 6000
$ws.Cells(136, 13).Value = 5750
$ws.Cells(136, 16).Value = 288
$ws.Cells(137, 4).Value = 44202
$ws.Cells(137, 10).Value = 720
$ws.Cells(137, 11).Value = 5000
$ws.Cells(137, 12).Value = 5500
$ws.Cells(137, 13).Value = 5250
$ws.Cells(137, 16).Value = 262
$ws.Cells(138, 4).Value = 44181
$ws.Cells(138, 10).Value = 400
$ws.Cells(139, 4).Value = 44369
$ws.Cells(139, 10).Value = 600
$ws.Cells(139, 11).Value = 5000
$ws.Cells(139, 12).Value = 5500
$ws.Cells(139, 13).Value = 5250
$ws.Cells(139, 16).Value = 262
$ws.Cells(140, 4).Value = 44307
$ws.Cells(140, 10).Value = 800
$ws.Cells(140, 11).Value = 5000
$ws.Cells(140, 12).Value = 5500
$ws.Cells(140, 13).Value = 5250
$ws.Cells(140, 16).Value = 262
$ws.Cells(141, 4).Value = 44161
$ws.Cells(141, 10).Value = 600
$ws.Cells(141, 11).Value = 4800
$ws.Cells(141, 12).Value = 5000
$ws.Cells(141, 13).Value = 4900
$ws.Cells(141, 16).Value = 245
$ws.Cells(142, 4).Value = 44195
$ws.Cells(142, 10).Value = 720
$ws.Cells(143, 4).Value = 44265
$ws.Cells(143, 11).Value = 5500
$ws.Cells(143, 12).Value = 6000
$ws.Cells(143, 13).Value = 5750
$ws.Cells(143, 16).Value = 288
$ws.Cells(144, 4).Value = 44159
$ws.Cells(144, 10).Value = 560
$ws.Cells(145, 4).Value = 44333
$ws.Cells(145, 10).Value = 660
$ws.Cells(146, 4).Value = 44201
$ws.Cells(146, 10).Value = 560
$ws.Cells(147, 4).Value = 44433
$ws.Cells(147, 10).Value = 800
$ws.Cells(147, 11).Value = 5000
$ws.Cells(147, 12).Value = 5500
$ws.Cells(147, 13).Value = 5250
$ws.Cells(147, 16).Value = 262
$ws.Cells(148, 4).Value = 44309
$ws.Cells(149, 4).Value = 44344
$ws.Cells(149, 10).Value = 800
$ws.Cells(150, 4).Value = 44319
$ws.Cells(150, 10).Value = 680
$ws.Cells(151, 4).Value = 44316
$ws.Cells(151, 10).Value = 800
$ws.Cells(152, 4).Value = 44253
$ws.Cells(152, 10).Value = 840
$ws.Cells(152, 11).Value = 5500
$ws.Cells(152, 12).Value = 6000
$ws.Cells(152, 13).Value = 5750
$ws.Cells(152, 16).Value = 288
$ws.Cells(153, 4).Value = 44281
$ws.Cells(153, 11).Value = 5000
$ws.Cells(153, 12).Value = 5500
$ws.Cells(153, 13).Value = 5250
$ws.Cells(153, 16).Value = 262
$ws.Cells(154, 4).Value = 44160
$ws.Cells(154, 10).Value = 720
$ws.Cells(154, 11).Value = 4800
$ws.Cells(154, 12).Value = 5000
$ws.Cells(154, 13).Value = 4900
$ws.Cells(154, 16).Value = 245
$ws.Cells(155, 4).Value = 44186
$ws.Cells(155, 11).Value = 5000
$ws.Cells(155, 12).Value = 5500
$ws.Cells(155, 13).Value = 5250
$ws.Cells(155, 16).Value = 262
$ws.Cells(156, 4).Value = 44211
$ws.Cells(156, 10).Value = 740
$ws.Cells(157, 4).Value = 44301
$ws.Cells(158, 4).Value = 44272
$ws.Cells(158, 10).Value = 800
$ws.Cells(158, 11).Value = 5500
$ws.Cells(158, 12).Value = 6000
$ws.Cells(158, 13).Value = 5750
$ws.Cells(158, 16).Value = 288
$ws.Cells(159, 4).Value = 44370
$ws.Cells(160, 4).Value = 44326
$ws.Cells(160, 10).Value = 700
$ws.Cells(160, 11).Value = 4800
$ws.Cells(160, 12).Value = 5000
$ws.Cells(160, 13).Value = 4900
$ws.Cells(160, 16).Value = 245
$ws.Cells(161, 4).Value = 44230
$ws.Cells(161, 10).Value = 720
$ws.Cells(162, 4).Value = 44382
$ws.Cells(162, 10).Value = 700
$ws.Cells(163, 4).Value = 44232
$ws.Cells(163, 10).Value = 760
$ws.Cells(163, 11).Value = 4800
$ws.Cells(163, 12).Value = 5000
$ws.Cells(163, 13).Value = 4900
$ws.Cells(163, 16).Value = 245
$ws.Cells(164, 4).Value = 44398
$ws.Cells(164, 10).Value = 800
$ws.Cells(164, 11).Value = 5000
$ws.Cells(164, 12).Value = 5500
$ws.Cells(164, 13).Value = 5250
$ws.Cells(164, 16).Value = 262
$ws.Cells(165, 4).Value = 44270
$ws.Cells(165, 10).Value = 600
$ws.Cells(165, 11).Value = 5500
$ws.Cells(165, 12).Value = 6000
$ws.Cells(165, 13).Value = 5750
$ws.Cells(165, 16).Value = 288
$ws.Cells(166, 4).Value = 44286
$ws.Cells(166, 10).Value = 800
$ws.Cells(167, 4).Value = 44217
$ws.Cells(167, 10).Value = 600
$ws.Cells(167, 12).Value = 5500
$ws.Cells(167, 13).Value = 5250
$ws.Cells(167, 16).Value = 262
$ws.Cells(168, 4).Value = 44244
$ws.Cells(168, 10).Value = 800
$ws.Cells(168, 11).Value = 5500
$ws.Cells(168, 12).Value = 6000
$ws.Cells(168, 13).Value = 5750
$ws.Cells(168, 16).Value = 288
$ws.Cells(169, 4).Value = 44343
$ws.Cells(169, 11).Value = 4800
$ws.Cells(169, 12).Value = 5000
$ws.Cells(169, 13).Value = 4900
$ws.Cells(169, 16).Value = 245
$ws.Cells(170, 4).Value = 44223
$ws.Cells(170, 10).Value = 740
$ws.Cells(171, 4).Value = 44376
$ws.Cells(172, 4).Value = 44208
$ws.Cells(172, 10).Value = 560
$ws.Cells(172, 11).Value = 5000
$ws.Cells(172, 12).Value = 6000
$ws.Cells(172, 13).Value = 5500
$ws.Cells(172, 16).Value = 275
$ws.Cells(173, 1).Value = 8
$ws.Cells(173, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells(173, 3).Value = "Coquimbo"
$ws.Cells(173, 4).Value = 44179
$ws.Cells(173, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells(173, 5).Value = 4
$ws.Cells(173, 6).Value = 100114013
$ws.Cells(173, 7).Value = "Zanahoria"
$ws.Cells(173, 8).Value = "Sin especificar"
$ws.Cells(173, 9).Value = "Primera"
$ws.Cells(173, 10).Value = 760
$ws.Cells(173, 11).Value = 5000
$ws.Cells(173, 12).Value = 5500
$ws.Cells(173, 13).Value = 5250
$ws.Cells(173, 14).Value = "`$/saco 20 kilos"
$ws.Cells(173, 15).Value = "Provincia del Elquí"
$ws.Cells(173, 16).Value = 262
$ws.Cells(173, 17).Value = 20
$ws.Cells(173, 18).Value = "Hortaliza"
$ws.Cells(174, 1).Value = 8
$ws.Cells(174, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells(174, 3).Value = "Coquimbo"
$ws.Cells(174, 4).Value = 44284
$ws.Cells(174, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells(174, 5).Value = 4
$ws.Cells(174, 6).Value = 100114013
$ws.Cells(174, 7).Value = "Zanahoria"
$ws.Cells(174, 8).Value = "Sin especificar"
$ws.Cells(174, 9).Value = "Primera"
$ws.Cells(174, 10).Value = 700
$ws.Cells(174, 11).Value = 5000
$ws.Cells(174, 12).Value = 5500
$ws.Cells(174, 13).Value = 5250
$ws.Cells(174, 14).Value = "`$/saco 20 kilos"
$ws.Cells(174, 15).Value = "Provincia del Elquí"
$ws.Cells(174, 16).Value = 262
$ws.Cells(174, 17).Value = 20
$ws.Cells(174, 18).Value = "Hortaliza"
$ws.Cells(175, 1).Value = 8
$ws.Cells(175, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells(175, 3).Value = "Coquimbo"
$ws.Cells(175, 4).Value = 44334
$ws.Cells(175, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells(175, 5).Value = 4
$ws.Cells(175, 6).Value = 100114013
$ws.Cells(175, 7).Value = "Zanahoria"
$ws.Cells(175, 8).Value = "Sin especificar"
$ws.Cells(175, 9).Value = "Primera"
$ws.Cells(175, 10).Value = 600
$ws.Cells(175, 11).Value = 5000
$ws.Cells(175, 12).Value = 5500
$ws.Cells(175, 13).Value = 5250
$ws.Cells(175, 14).Value = "`$/saco 20 kilos"
$ws.Cells(175, 15).Value = "Provincia del Elquí"
$ws.Cells(175, 16).Value = 262
$ws.Cells(175, 17).Value = 20
$ws.Cells(175, 18).Value = "Hortaliza"
$ws.Cells(176, 1).Value = 8
$ws.Cells(176, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells(176, 3).Value = "Coquimbo"
$ws.Cells(176, 4).Value = 44168
$ws.Cells(176, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells(176, 5).Value = 4
$ws.Cells(176, 6).Value = 100114013
$ws.Cells(176, 7).Value = "Zanahoria"
$ws.Cells(176, 8).Value = "Sin especificar"
$ws.Cells(176, 9).Value = "Primera"
$ws.Cells(176, 10).Value = 600
$ws.Cells(176, 11).Value = 5000
$ws.Cells(176, 12).Value = 5500
$ws.Cells(176, 13).Value = 5250
$ws.Cells(176, 14).Value = "`$/saco 20 kilos"
$ws.Cells(176, 15).Value = "Provincia del Elquí"
$ws.Cells(176, 16).Value = 262
$ws.Cells(176, 17).Value = 20
$ws.Cells(176, 18).Value = "Hortaliza"
$ws.Cells(177, 1).Value = 8
$ws.Cells(177, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells(177, 3).Value = "Coquimbo"
$ws.Cells(177, 4).Value = 44418
$ws.Cells(177, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells(177, 5).Value = 4
$ws.Cells(177, 6).Value = 100114013
$ws.Cells(177, 7).Value = "Zanahoria"
$ws.Cells(177, 8).Value = "Sin especificar"
$ws.Cells(177, 9).Value = "Primera"
$ws.Cells(177, 10).Value = 600
$ws.Cells(177, 11).Value = 5000
$ws.Cells(177, 12).Value = 5500
$ws.Cells(177, 13).Value = 5250
$ws.Cells(177, 14).Value = "`$/saco 20 kilos"
$ws.Cells(177, 15).Value = "Provincia del Elquí"
$ws.Cells(177, 16).Value = 262
$ws.Cells(177, 17).Value = 20
$ws.Cells(177, 18).Value = "Hortaliza"
$ws.Cells(178, 1).Value = 8
$ws.Cells(178, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells(178, 3).Value = "Coquimbo"
$ws.Cells(178, 4).Value = 44432
$ws.Cells(178, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells(178, 5).Value = 4
$ws.Cells(178, 6).Value = 100114013
$ws.Cells(178, 7).Value = "Zanahoria"
$ws.Cells(178, 8).Value = "Sin especificar"
$ws.Cells(178, 9).Value = "Primera"
$ws.Cells(178, 10).Value = 660
$ws.Cells(178, 11).Value = 5000
$ws.Cells(178, 12).Value = 5500
$ws.Cells(178, 13).Value = 5250
$ws.Cells(178, 14).Value = "`$/saco 20 kilos"
$ws.Cells(178, 15).Value = "Provincia del Elquí"
$ws.Cells(178, 16).Value = 262
$ws.Cells(178, 17).Value = 20
$ws.Cells(178, 18).Value = "Hortaliza"
